$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.771.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.62%  "
$ws.Range("D3").Value = "'3.011.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'565.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.86%  "
$ws.Range("D6").Value = "'140.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.43%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").Value = "'3.002.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").Value = "'0.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.81%  "
$ws.Range("D11").Value = "'5.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.46%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "'0.0000232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("D14").Value = "'33.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'3.518.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("D18").Value = "'3.015.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("D19").Value = "'59.770.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").Value = "'438.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.75%  "
$ws.Range("D21").Value = "'13.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("D22").Value = "'0.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("D23").Value = "'7.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "'13.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("D25").Value = "'80.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'2.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.95%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'2.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("D30").Value = "'7.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.38%  "
$ws.Range("D31").Value = "'6.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.47%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.107"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.42%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'25.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").Value = "'0.0₃0787"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.57%  "
$ws.Range("E35").Value = "  +6.82%  "
$ws.Range("D36").Value = "'5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.47%  "
$ws.Range("D37").Value = "'2.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "'49.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "'8.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.44%  "
$ws.Range("D41").Value = "'403.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.95%  "
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("D43").Value = "'2.767.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("D44").Value = "'0.107"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "'0.253"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.38%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'123.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "'2.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.96%  "
$ws.Range("D49").Value = "'0.110"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").Value = "'34.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +19.56%  "
$ws.Range("D51").Value = "'23.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.33%  "
